$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130 (pushes existing rows 130.. down by one),
# matching the weekly data-entry pattern used throughout this sheet.
$ws.Rows("130:130").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A130").Value = 3
$ws.Range("B130").Value = "Femacal de La Calera"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44566
$ws.Range("E130").Value = 5
$ws.Range("F130").Value = 100112012
$ws.Range("G130").Value = "Espinaca"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 170
$ws.Range("K130").Value = 3500
$ws.Range("L130").Value = 4000
$ws.Range("M130").Value = 3765
$ws.Range("N130").Value = '$/docena de atados (3 kilos)'
$ws.Range("O130").Value = "Provincia de Quillota"
$ws.Range("P130").Value = 1255
$ws.Range("Q130").Value = 3
$ws.Range("R130").Value = "Hortaliza"
